$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (E1) entirely, shifting F:K left to E:J.
$ws.Range("E1").EntireColumn.Delete()
